$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22 ("More Enhancements") so that it
# becomes row 23, and populate the new row 22 with the new task.
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "Weather Prediction Enhancements"

# Give B22 the same (date) number format as the other date cells in column B,
# but leave it without a value, matching the rest of the timeline structure.
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A23").Select()
